# Atualização do BurnDown e Valor Agregado
# Task rows 56 and 58 ("A.A.A" / "A.A.A") now have actual-cost hours
# reported for Sprint 6, and are marked as completed ("S") for Sprint 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 56: add Actual Cost (B56) and mark Sprint 6 complete (I56) ---
$ws.Range("B56").Value = 5
$ws.Range("I56").Value = "S"

# --- Row 58: add Actual Cost (B58) and mark Sprint 6 complete (I58) ---
$ws.Range("B58").Value = 1
$ws.Range("I58").Value = "S"

# Update the view state to reflect where the author ended up after the edit
$ws.Range("B86").Select()
